$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the source data.
# Deleting from the bottom up (row 28 before row 26) keeps the earlier,
# still-pending row number ("RM 232" at row 26) valid.
$ws.Rows(28).Delete()   # "SC 92" row removed
$ws.Rows(26).Delete()   # "RM 232" row removed

# Individual value corrections (imputed / cleared cells) after the shift.
$ws.Range("E3").Value = -5.7
$ws.Range("F4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F17").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").Value = ""

# After the two row deletions, the "SC 193" row (now row 32) needs its
# previously-missing D column (error) value restored.
$ws.Range("E32").Value = -6.4
